$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlineLogo($headerFooter, $newName) {
    if ($headerFooter.Exists -and $headerFooter.Range.InlineShapes.Count -ge 1) {
        $inlineShape = $headerFooter.Range.InlineShapes(1)
        # Renaming InlineShapes directly is unreliable for header/footer
        # stories in this host, so round-trip through a floating Shape
        # (where the rename is applied) and back to an inline shape.
        $floatingShape = $inlineShape.ConvertToShape()
        $floatingShape.Name = $newName
        $floatingShape.ConvertToInlineShape()
    }
}

# Footer (default / file footer2.xml, docPr id="2"): PearsonLogo image2.png -> image1.png
Rename-InlineLogo $sec.Footers(1) "image1.png"

# Footer (first page / file footer1.xml, docPr id="3"): PearsonLogo image2.png -> image1.png
Rename-InlineLogo $sec.Footers(2) "image1.png"

# Header (first page / file header1.xml, docPr id="1"): BTec_Logo-Orange image1.jpg -> image2.jpg
Rename-InlineLogo $sec.Headers(2) "image2.jpg"

Write-Host "renamed logos"
